$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header labels in row 1 (columns N, O, P) with the new
# parenthetical annotations, e.g. "Общая сумма (3)" -> "Общая сумма (0-1)".
$ws.Range("N1").Value = "Общая сумма (0-1)"
$ws.Range("O1").Value = "Общие комиссионные (1)"
$ws.Range("P1").Value = "Всего (0+1)"

# Move the active selection to P2 (was N8).
$ws.Range("P2").Select()
